$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Recent events further sparked interest ... whether classical
#    statistical methods" -> "Recent events have further simulated interest
#    ... whether it is classic statistical methods"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Recent events further sparked interest in learning time series analysis and forecasting techniques,  whether classical statistical methods",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Recent events have further simulated interest in learning time series analysis and prediction techniques,  whether it is classic statistical methods",
    2)

# ---------------------------------------------------------------------------
# 2) "predicting presidential election outcome," -> "... results,"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "predicting presidential election outcome,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "predicting presidential election results,",
    2)

# ---------------------------------------------------------------------------
# 3) " Marketing planning , " -> drop the space before the comma:
#    " Marketing planning, "
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Marketing planning , Inventory",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Marketing planning, Inventory",
    2)

# ---------------------------------------------------------------------------
# 4) Turn the trailing empty "L-Bullets" paragraph (the one with the
#    numPr/ind direct-formatting override) into three new "P-Regular"
#    paragraphs of body text.
# ---------------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "L - Bullets" -and $p.Range.Text.Trim() -eq "") {
        $target = $p
        break
    }
}

$target.Style = "P-Regular"
$target.Range.InsertAfter("When working with time series data we usually have two goals: time series analysis (descriptive), and time series forecasting (predictive). In time series analysis we strive to extract a better understanding and deeper intuition into the underlying phenomenon captured in our data using statistical methods. In time series forecasting, we aim to create a predictive model that extends from our data into the unforeseeable future and hence the term forecasting.  ")

$target.Range.InsertParagraphAfter()
$para2 = $target.Next()
$para2.Style = "P-Regular"
$para2.Range.InsertAfter("Time series data differs from the typical data used in machine learning in the classical sense due to the dependence on time, serial correlation, and dependence. ")

$para2.Range.InsertParagraphAfter()
$para3 = $para2.Next()
$para3.Style = "P-Regular"
$para3.Range.InsertAfter("In this book, we cover a variety of recipes for both analysis and forecasting of time series data and pragmatic approach to handling the complex nature of time series data. ")
